# Remove the (now-unused) leading column A (previously blank/placeholder).
# This shifts B:D left into A:C, so:
#   B1/C1/D1 (headers)      -> A1/B1/C1   (keeps the header style s="1")
#   old A2 (value 0, s="1") is deleted entirely
#   B2/C2/D2 (values 1,2,3) -> A2/B2/C2   (no special style, like before)
# The sheet dimension shrinks from A1:D2 to A1:C2 automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).Delete()
